$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: split the single run in the "chosen topic" paragraph into
# three runs, wrapping "Learning ." with a gramStart/gramEnd proofErr
# pair (content/text is unchanged, only run boundaries change).
# -----------------------------------------------------------------
$target1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("The chosen topic is Big Data and Neural Networks")) {
        $target1 = $p
        break
    }
}
if ($target1 -eq $null) {
    throw "Could not find 'The chosen topic' paragraph"
}

$r1 = $target1.Range
# Exclude the trailing paragraph mark so only the run content is replaced.
$r1 = $d.Range($r1.Start, $r1.End - 1)

$xml1 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t xml:space="preserve">The chosen topic is Big Data and Neural Networks, with NN being considered a type of Machine Learning (ML) process known as Deep </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Learning .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> The field of Big Data is constantly growing and encompasses a need for efficient data management and processing tools. Two well-known tools for handling and analyzing large datasets are Relational Database Management Systems (RDBMS) and Hadoop. However, the rampant advancement of Machine Learning and Neural Networks, the integration of these data management tools with advanced analytics technologies is the focus of this paper.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$null = $r1.InsertXML($xml1)

# -----------------------------------------------------------------
# Edit 2: replace the "The current state of RDBMS ..." paragraph and
# the following empty paragraph with a whole new block of content:
#   - the (untouched) RDBMS paragraph, but now left-aligned
#   - a blank separator paragraph
#   - a new Hadoop paragraph
#   - a blank separator paragraph
#   - a new Machine Learning APIs paragraph
#   - a blank separator paragraph
#   - a new Deep Learning Neural Networks paragraph (reuses the old
#     "both" justified / first-line-indent pPr)
#   - a trailing blank separator paragraph
# -----------------------------------------------------------------
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("The current state of RDBMS has evolved substantially")) {
        $startIndex = $i
        break
    }
}
if ($startIndex -eq -1) {
    throw "Could not find 'The current state of RDBMS' paragraph"
}
$startPara = $d.Paragraphs.Item($startIndex)
$endPara = $d.Paragraphs.Item($startIndex + 1)

$full = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml2 = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:ind w:firstLine="288"/><w:jc w:val="left"/></w:pPr><w:r><w:t>The current state of RDBMS has evolved substantially with enhancements in storage, speed, and scalability by using cloud-based solutions</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>The future holds a shift for RDBMS transitioning to a NoSQL database</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>To understand why NoSQL is taking over RDBMS, it is crucial to talk about: Schemas, where NoSQL uses dynamic instead of static schemas; the type of data to be stored, with NoSQL databases offering advantages for hierarchical data storage due to their flexible data models and scalability, while RDBMS are not that flexible; scalability, with NoSQL depending on horizontal scalability and RDBMS on vertical scalability; and other points where NoSQL surpasses RDBMS, including data warehouse, complexity, cloud, and big data handling, and output performance</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:ind w:firstLine="288"/><w:jc w:val="left"/></w:pPr></w:p>
<w:p><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:t>In the era of Big Data, where we sometimes run out of storage and face difficulties on a single host due to the volume of data,</w:t></w:r><w:r><w:t xml:space="preserve"> Hadoop came </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>into</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the scene to tackle this by offering computational capabilities over huge amounts of data</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">The present and future look bright for Hadoop, as some of the major Big Data companies, such as </w:t></w:r><w:r><w:t xml:space="preserve">Google, Facebook, eBay, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Twitter</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t xml:space="preserve">Spotify, </w:t></w:r><w:r><w:t>rely on this technology</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:jc w:val="left"/></w:pPr></w:p>
<w:p><w:pPr><w:ind w:firstLine="288"/><w:jc w:val="left"/></w:pPr><w:r><w:t>Machine Learning APIs have helped developers integrating data flows into complex algorithms without requiring deep expertise</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">These APIs were once primarily used for basic tasks like picture and speech </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>recognition, but they have since grown to include a variety of machine learning activities, such as predictive analytics and natural language processing</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Today, they are essential to leading tech companies and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>cutting edge</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> industries, being ML accessible and customizable than in the past</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Looking</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> at the future ML APIs will remain as they have simplified model development across diverse environments</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
<w:p><w:pPr><w:jc w:val="left"/></w:pPr></w:p>
<w:p><w:pPr><w:ind w:firstLine="288"/><w:jc w:val="both"/></w:pPr><w:r><w:t>The development of Deep Learning Neural Networks (DLNNs) traces back to the 1950s. The method has improved since the introduction of Convolutional Neural Networks (CNNs) by LeCun et al. in the late 1980s, which showed how good deep architectures are for image processing</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">In 2012 this field achieved another milestone with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AlexNet's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> success in the ImageNet challenge, showcasing DLNNs' potential in image recognition</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>The introduction of Transformer models by Vaswani et al. in 2017 marked another significant advancement, in natural language processing. Neural Networks have the future guaranteed as all points mentioned earlier RDMS and Hadoop store and process data for Neural Networks models.</w:t></w:r></w:p>
<w:p><w:pPr><w:jc w:val="left"/></w:pPr></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$null = $full.InsertXML($xml2)

Write-Output "Edits applied."
